$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: update existing metric values (MSE, R2, MAE) ---
$ws.Range("B2").Value = 0.06146626137578159
$ws.Range("C2").Value = 0.9994169649648684
$ws.Range("D2").Value = 0.1798812967085154

# --- New column F: "Modelo" ---
# F1 header: copy the formatting from the adjacent header cell (E1, "Tipo")
# so it gets the same bold/bordered/centered header style, then set its text.
$ws.Range("F1").Value = "Modelo"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

# F2 data cell: plain value, same as the rest of row 2.
$ws.Range("F2").Value = "Pipeline(steps=[('model', RandomForestRegressor(max_depth=5, n_estimators=50))])"
